$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 110 (new weekly price entry), pushing the
# existing rows 110-115 down to 111-116. Excel copies the formatting
# (incl. the date style on column D) from the row above automatically.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly entry.
$ws.Cells.Item(110, 1).Value = 10
$ws.Cells.Item(110, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(110, 3).Value = "La Araucanía"
$ws.Cells.Item(110, 4).Value = 44461
$ws.Cells.Item(110, 5).Value = 9
$ws.Cells.Item(110, 6).Value = 100112013
$ws.Cells.Item(110, 7).Value = "Alcachofa"
$ws.Cells.Item(110, 8).Value = "Madrigal"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 100
$ws.Cells.Item(110, 11).Value = 12000
$ws.Cells.Item(110, 12).Value = 12000
$ws.Cells.Item(110, 13).Value = 12000
$ws.Cells.Item(110, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(110, 15).Value = "Región Metropolitana"
$ws.Cells.Item(110, 16).Value = 300
$ws.Cells.Item(110, 17).Value = 40
$ws.Cells.Item(110, 18).Value = "Hortaliza"
